# Add a placeholder row for a real (not-yet-sequenced) metagenome sample
# "ZymoFecal" (row 8, column A already has the sample name).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Date of sequencing (B8) - use same date format as the existing
# " date of basecalling" column (G2/G3) by copying its format over.
$ws.Range("B8").Value = 45566
$ws.Range("G2").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Sequencing kit / flowcell type
$ws.Range("C8").Value = "LSK114"
$ws.Range("D8").Value = "FLO-PRO114M"

# Sample rate is written before the flowcell ID so that the shared-string
# table ends up with the same ordering as the source workbook.
$ws.Range("F8").Value = "5 khz"
$ws.Range("E8").Value = "PAW77640"

$ws.Range("G8").Select()
